# paises.xlsx update: refresh COVID-19 country stats and re-sort by total cases,
# plus bump the "Datos actualizados" timestamp (commit: "Update countries & provincias Spain").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: refresh timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Junio de 2020 a las 19:04"

# Rows 4-219: country name (col A) + Casos totales/Nuevos casos/Casos activos/
# Recuperados/Casos criticos/Muertes hoy/Muertes (cols B-H).
# Only rows whose country or figures actually changed are touched below;
# everything else keeps its original value.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 2343104
$ws.Cells.Item(4, 3).Value = 12526
$ws.Cells.Item(4, 4).Value = 974757
$ws.Cells.Item(4, 5).Value = 1246278
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 89
$ws.Cells.Item(4, 8).Value = 122069

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 1073376
$ws.Cells.Item(5, 3).Value = 3237
$ws.Cells.Item(5, 4).Value = 543186
$ws.Cells.Item(5, 5).Value = 480008
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 124
$ws.Cells.Item(5, 8).Value = 50182

# Row 7: India
$ws.Cells.Item(7, 1).Value = "India"
$ws.Cells.Item(7, 2).Value = 422526
$ws.Cells.Item(7, 3).Value = 10799
$ws.Cells.Item(7, 4).Value = 235328
$ws.Cells.Item(7, 5).Value = 173689
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 232
$ws.Cells.Item(7, 8).Value = 13509

# Row 9: España
$ws.Cells.Item(9, 1).Value = "España"
$ws.Cells.Item(9, 2).Value = 293352
$ws.Cells.Item(9, 3).Value = 334
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 28323

# Row 11: Italia
$ws.Cells.Item(11, 1).Value = "Italia"
$ws.Cells.Item(11, 2).Value = 238499
$ws.Cells.Item(11, 3).Value = 224
$ws.Cells.Item(11, 4).Value = 182893
$ws.Cells.Item(11, 5).Value = 20972
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 24
$ws.Cells.Item(11, 8).Value = 34634

# Row 30: Ecuador
$ws.Cells.Item(30, 1).Value = "Ecuador"
$ws.Cells.Item(30, 2).Value = 50640
$ws.Cells.Item(30, 3).Value = 909
$ws.Cells.Item(30, 4).Value = 24991
$ws.Cells.Item(30, 5).Value = 21426
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 67
$ws.Cells.Item(30, 8).Value = 4223

# Row 34: Singapur
$ws.Cells.Item(34, 1).Value = "Singapur"
$ws.Cells.Item(34, 2).Value = 42095
$ws.Cells.Item(34, 3).Value = 262
$ws.Cells.Item(34, 4).Value = 34942
$ws.Cells.Item(34, 5).Value = 7127
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 26

# Row 41: Irak
$ws.Cells.Item(41, 1).Value = "Irak"
$ws.Cells.Item(41, 2).Value = 30868
$ws.Cells.Item(41, 3).Value = 1646
$ws.Cells.Item(41, 4).Value = 13935
$ws.Cells.Item(41, 5).Value = 15833
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 87
$ws.Cells.Item(41, 8).Value = 1100

# Row 42: Filipinas
$ws.Cells.Item(42, 1).Value = "Filipinas"
$ws.Cells.Item(42, 2).Value = 30052
$ws.Cells.Item(42, 3).Value = 652
$ws.Cells.Item(42, 4).Value = 7893
$ws.Cells.Item(42, 5).Value = 20990
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 19
$ws.Cells.Item(42, 8).Value = 1169

# Row 43: Oman
$ws.Cells.Item(43, 1).Value = "Oman"
$ws.Cells.Item(43, 2).Value = 29471
$ws.Cells.Item(43, 3).Value = 905
$ws.Cells.Item(43, 4).Value = 15552
$ws.Cells.Item(43, 5).Value = 13788
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 3
$ws.Cells.Item(43, 8).Value = 131

# Row 46: Irlanda
$ws.Cells.Item(46, 1).Value = "Irlanda"
$ws.Cells.Item(46, 2).Value = 25379
$ws.Cells.Item(46, 3).Value = 5
$ws.Cells.Item(46, 4).Value = 22698
$ws.Cells.Item(46, 5).Value = 966
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 1715

# Row 56: Kazajistan
$ws.Cells.Item(56, 1).Value = "Kazajistan"
$ws.Cells.Item(56, 2).Value = 17225
$ws.Cells.Item(56, 3).Value = 446
$ws.Cells.Item(56, 4).Value = 10897
$ws.Cells.Item(56, 5).Value = 6208
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 2
$ws.Cells.Item(56, 8).Value = 120

# Row 61: Azerbaiyan
$ws.Cells.Item(61, 1).Value = "Azerbaiyan"
$ws.Cells.Item(61, 2).Value = 12729
$ws.Cells.Item(61, 3).Value = 491
$ws.Cells.Item(61, 4).Value = 6799
$ws.Cells.Item(61, 5).Value = 5776
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 6
$ws.Cells.Item(61, 8).Value = 154

# Row 62: Corea del Sur
$ws.Cells.Item(62, 1).Value = "Corea del Sur"
$ws.Cells.Item(62, 2).Value = 12421
$ws.Cells.Item(62, 3).Value = 48
$ws.Cells.Item(62, 4).Value = 10868
$ws.Cells.Item(62, 5).Value = 1273
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 280

# Row 63: Dinamarca
$ws.Cells.Item(63, 1).Value = "Dinamarca"
$ws.Cells.Item(63, 2).Value = 12391
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 11282
$ws.Cells.Item(63, 5).Value = 509
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 600

# Row 64: Honduras
$ws.Cells.Item(64, 1).Value = "Honduras"
$ws.Cells.Item(64, 2).Value = 12306
$ws.Cells.Item(64, 3).Value = 1048
$ws.Cells.Item(64, 4).Value = 1275
$ws.Cells.Item(64, 5).Value = 10673
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 9
$ws.Cells.Item(64, 8).Value = 358

# Row 65: Argelia
$ws.Cells.Item(65, 1).Value = "Argelia"
$ws.Cells.Item(65, 2).Value = 11771
$ws.Cells.Item(65, 3).Value = 140
$ws.Cells.Item(65, 4).Value = 8422
$ws.Cells.Item(65, 5).Value = 2504
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 8
$ws.Cells.Item(65, 8).Value = 845

# Row 83: Kenia
$ws.Cells.Item(83, 1).Value = "Kenia"
$ws.Cells.Item(83, 2).Value = 4738
$ws.Cells.Item(83, 3).Value = 260
$ws.Cells.Item(83, 4).Value = 1607
$ws.Cells.Item(83, 5).Value = 3008
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 123

# Row 84: El Salvador
$ws.Cells.Item(84, 1).Value = "El Salvador"
$ws.Cells.Item(84, 2).Value = 4626
$ws.Cells.Item(84, 3).Value = 151
$ws.Cells.Item(84, 4).Value = 2535
$ws.Cells.Item(84, 5).Value = 1993
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 5
$ws.Cells.Item(84, 8).Value = 98

# Row 85: Republica de Yibuti
$ws.Cells.Item(85, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(85, 2).Value = 4582
$ws.Cells.Item(85, 3).Value = 17
$ws.Cells.Item(85, 4).Value = 3859
$ws.Cells.Item(85, 5).Value = 678
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 45

# Row 86: Etiopia
$ws.Cells.Item(86, 1).Value = "Etiopia"
$ws.Cells.Item(86, 2).Value = 4532
$ws.Cells.Item(86, 3).Value = 63
$ws.Cells.Item(86, 4).Value = 1213
$ws.Cells.Item(86, 5).Value = 3245
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 74

# Row 98: Somalia
$ws.Cells.Item(98, 1).Value = "Somalia"
$ws.Cells.Item(98, 2).Value = 2779
$ws.Cells.Item(98, 3).Value = 24
$ws.Cells.Item(98, 4).Value = 782
$ws.Cells.Item(98, 5).Value = 1907
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 90

# Row 99: Guayana Francesa
$ws.Cells.Item(99, 1).Value = "Guayana Francesa"
$ws.Cells.Item(99, 2).Value = 2441
$ws.Cells.Item(99, 3).Value = 278
$ws.Cells.Item(99, 4).Value = 930
$ws.Cells.Item(99, 5).Value = 1505
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 6

# Row 100: Mayotte
$ws.Cells.Item(100, 1).Value = "Mayotte"
$ws.Cells.Item(100, 2).Value = 2404
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 2066
$ws.Cells.Item(100, 5).Value = 307
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 31

# Row 101: Croacia
$ws.Cells.Item(101, 1).Value = "Croacia"
$ws.Cells.Item(101, 2).Value = 2317
$ws.Cells.Item(101, 3).Value = 18
$ws.Cells.Item(101, 4).Value = 2142
$ws.Cells.Item(101, 5).Value = 68
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 107

# Row 102: Cuba
$ws.Cells.Item(102, 1).Value = "Cuba"
$ws.Cells.Item(102, 2).Value = 2312
$ws.Cells.Item(102, 3).Value = 3
$ws.Cells.Item(102, 4).Value = 2103
$ws.Cells.Item(102, 5).Value = 124
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 85

# Row 103: Maldivas
$ws.Cells.Item(103, 1).Value = "Maldivas"
$ws.Cells.Item(103, 2).Value = 2187
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 1788
$ws.Cells.Item(103, 5).Value = 391
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 8

# Row 122: Sierra Leona
$ws.Cells.Item(122, 1).Value = "Sierra Leona"
$ws.Cells.Item(122, 2).Value = 1327
$ws.Cells.Item(122, 3).Value = 18
$ws.Cells.Item(122, 4).Value = 788
$ws.Cells.Item(122, 5).Value = 484
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 2
$ws.Cells.Item(122, 8).Value = 55

# Row 128: Republica de Chipre
$ws.Cells.Item(128, 1).Value = "Republica de Chipre"
$ws.Cells.Item(128, 2).Value = 986
$ws.Cells.Item(128, 3).Value = 1
$ws.Cells.Item(128, 4).Value = 824
$ws.Cells.Item(128, 5).Value = 143
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 19

# Row 129: Yemen
$ws.Cells.Item(129, 1).Value = "Yemen"
$ws.Cells.Item(129, 2).Value = 941
$ws.Cells.Item(129, 3).Value = 19
$ws.Cells.Item(129, 4).Value = 347
$ws.Cells.Item(129, 5).Value = 338
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 2
$ws.Cells.Item(129, 8).Value = 256

# Row 135: Republica del Chad
$ws.Cells.Item(135, 1).Value = "Republica del Chad"
$ws.Cells.Item(135, 2).Value = 858
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 752
$ws.Cells.Item(135, 5).Value = 32
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 74

# Row 202: Dominica
$ws.Cells.Item(202, 1).Value = "Dominica"
$ws.Cells.Item(202, 2).Value = 18
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 18
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203: Fiyi
$ws.Cells.Item(203, 1).Value = "Fiyi"
$ws.Cells.Item(203, 2).Value = 18
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 18
$ws.Cells.Item(203, 5).Value = 0
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

# Row 206: Islas Turcas y Caicos
$ws.Cells.Item(206, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(206, 2).Value = 14
$ws.Cells.Item(206, 3).Value = 2
$ws.Cells.Item(206, 4).Value = 11
$ws.Cells.Item(206, 5).Value = 2
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 1

# Row 207: Islas Malvinas
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 208: Groenlandia
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 2).Value = 13
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 13
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

# Row 209: Santa Sede
$ws.Cells.Item(209, 1).Value = "Santa Sede"
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 213: Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Papua Nueva Guinea
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0
